$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dArr = New-Object 'object[,]' 36,1
$dArr[0,0] = 44441
$dArr[1,0] = 44455
$dArr[2,0] = 44466
$dArr[3,0] = 44497
$dArr[4,0] = 44497
$dArr[5,0] = 44496
$dArr[6,0] = 44438
$dArr[7,0] = 44448
$dArr[8,0] = 44369
$dArr[9,0] = 44369
$dArr[10,0] = 44186
$dArr[11,0] = 44179
$dArr[12,0] = 44424
$dArr[13,0] = 44512
$dArr[14,0] = 44294
$dArr[15,0] = 44315
$dArr[16,0] = 44315
$dArr[17,0] = 44316
$dArr[18,0] = 44364
$dArr[19,0] = 44525
$dArr[20,0] = 44508
$dArr[21,0] = 44532
$dArr[22,0] = 44498
$dArr[23,0] = 44425
$dArr[24,0] = 44161
$dArr[25,0] = 44341
$dArr[26,0] = 44389
$dArr[27,0] = 44321
$dArr[28,0] = 44385
$dArr[29,0] = 44529
$dArr[30,0] = 44175
$dArr[31,0] = 44188
$dArr[32,0] = 44340
$dArr[33,0] = 44329
$dArr[34,0] = 44452
$dArr[35,0] = 44511
$ws.Range("D2:D37").Value = $dArr

$jqArr = New-Object 'object[,]' 36,8
$jqArr[0,0] = 40
$jqArr[0,1] = 20000
$jqArr[0,2] = 20000
$jqArr[0,3] = 20000
$jqArr[0,4] = '$/malla 20 kilos'
$jqArr[0,5] = 'Perú'
$jqArr[0,6] = 1000
$jqArr[0,7] = 20
$jqArr[1,0] = 30
$jqArr[1,1] = 20000
$jqArr[1,2] = 20000
$jqArr[1,3] = 20000
$jqArr[1,4] = '$/malla 20 kilos'
$jqArr[1,5] = 'Perú'
$jqArr[1,6] = 1000
$jqArr[1,7] = 20
$jqArr[2,0] = 20
$jqArr[2,1] = 25000
$jqArr[2,2] = 25000
$jqArr[2,3] = 25000
$jqArr[2,4] = '$/caja 15 kilos granel'
$jqArr[2,5] = 'Perú'
$jqArr[2,6] = 1667
$jqArr[2,7] = 15
$jqArr[3,0] = 30
$jqArr[3,1] = 20000
$jqArr[3,2] = 20000
$jqArr[3,3] = 20000
$jqArr[3,4] = '$/caja 15 kilos granel'
$jqArr[3,5] = 'Perú'
$jqArr[3,6] = 1333
$jqArr[3,7] = 15
$jqArr[4,0] = 40
$jqArr[4,1] = 20000
$jqArr[4,2] = 20000
$jqArr[4,3] = 20000
$jqArr[4,4] = '$/malla 20 kilos'
$jqArr[4,5] = 'Perú'
$jqArr[4,6] = 1000
$jqArr[4,7] = 20
$jqArr[5,0] = 30
$jqArr[5,1] = 20000
$jqArr[5,2] = 20000
$jqArr[5,3] = 20000
$jqArr[5,4] = '$/malla 20 kilos'
$jqArr[5,5] = 'Perú'
$jqArr[5,6] = 1000
$jqArr[5,7] = 20
$jqArr[6,0] = 40
$jqArr[6,1] = 20000
$jqArr[6,2] = 20000
$jqArr[6,3] = 20000
$jqArr[6,4] = '$/caja 15 kilos granel'
$jqArr[6,5] = 'Región de Arica y Parinacota'
$jqArr[6,6] = 1333
$jqArr[6,7] = 15
$jqArr[7,0] = 45
$jqArr[7,1] = 20000
$jqArr[7,2] = 20000
$jqArr[7,3] = 20000
$jqArr[7,4] = '$/malla 20 kilos'
$jqArr[7,5] = 'Perú'
$jqArr[7,6] = 1000
$jqArr[7,7] = 20
$jqArr[8,0] = 20
$jqArr[8,1] = 20000
$jqArr[8,2] = 20000
$jqArr[8,3] = 20000
$jqArr[8,4] = '$/caja 15 kilos granel'
$jqArr[8,5] = 'Región de Arica y Parinacota'
$jqArr[8,6] = 1333
$jqArr[8,7] = 15
$jqArr[9,0] = 20
$jqArr[9,1] = 20000
$jqArr[9,2] = 20000
$jqArr[9,3] = 20000
$jqArr[9,4] = '$/malla 20 kilos'
$jqArr[9,5] = 'Región de Arica y Parinacota'
$jqArr[9,6] = 1000
$jqArr[9,7] = 20
$jqArr[10,0] = 20
$jqArr[10,1] = 20000
$jqArr[10,2] = 20000
$jqArr[10,3] = 20000
$jqArr[10,4] = '$/caja 15 kilos granel'
$jqArr[10,5] = 'Región de Arica y Parinacota'
$jqArr[10,6] = 1333
$jqArr[10,7] = 15
$jqArr[11,0] = 20
$jqArr[11,1] = 20000
$jqArr[11,2] = 20000
$jqArr[11,3] = 20000
$jqArr[11,4] = '$/caja 15 kilos granel'
$jqArr[11,5] = 'Región de Arica y Parinacota'
$jqArr[11,6] = 1333
$jqArr[11,7] = 15
$jqArr[12,0] = 30
$jqArr[12,1] = 20000
$jqArr[12,2] = 20000
$jqArr[12,3] = 20000
$jqArr[12,4] = '$/caja 15 kilos granel'
$jqArr[12,5] = 'Región de Arica y Parinacota'
$jqArr[12,6] = 1333
$jqArr[12,7] = 15
$jqArr[13,0] = 30
$jqArr[13,1] = 20000
$jqArr[13,2] = 20000
$jqArr[13,3] = 20000
$jqArr[13,4] = '$/malla 20 kilos'
$jqArr[13,5] = 'Perú'
$jqArr[13,6] = 1000
$jqArr[13,7] = 20
$jqArr[14,0] = 5
$jqArr[14,1] = 20000
$jqArr[14,2] = 20000
$jqArr[14,3] = 20000
$jqArr[14,4] = '$/caja 15 kilos granel'
$jqArr[14,5] = 'Perú'
$jqArr[14,6] = 1333
$jqArr[14,7] = 15
$jqArr[15,0] = 30
$jqArr[15,1] = 20000
$jqArr[15,2] = 20000
$jqArr[15,3] = 20000
$jqArr[15,4] = '$/caja 15 kilos granel'
$jqArr[15,5] = 'Región de Arica y Parinacota'
$jqArr[15,6] = 1333
$jqArr[15,7] = 15
$jqArr[16,0] = 30
$jqArr[16,1] = 20000
$jqArr[16,2] = 20000
$jqArr[16,3] = 20000
$jqArr[16,4] = '$/malla 20 kilos'
$jqArr[16,5] = 'Región de Arica y Parinacota'
$jqArr[16,6] = 1000
$jqArr[16,7] = 20
$jqArr[17,0] = 20
$jqArr[17,1] = 20000
$jqArr[17,2] = 20000
$jqArr[17,3] = 20000
$jqArr[17,4] = '$/caja 15 kilos granel'
$jqArr[17,5] = 'Región de Arica y Parinacota'
$jqArr[17,6] = 1333
$jqArr[17,7] = 15
$jqArr[18,0] = 15
$jqArr[18,1] = 20000
$jqArr[18,2] = 20000
$jqArr[18,3] = 20000
$jqArr[18,4] = '$/caja 15 kilos granel'
$jqArr[18,5] = 'Perú'
$jqArr[18,6] = 1333
$jqArr[18,7] = 15
$jqArr[19,0] = 40
$jqArr[19,1] = 20000
$jqArr[19,2] = 20000
$jqArr[19,3] = 20000
$jqArr[19,4] = '$/caja 15 kilos granel'
$jqArr[19,5] = 'Perú'
$jqArr[19,6] = 1333
$jqArr[19,7] = 15
$jqArr[20,0] = 40
$jqArr[20,1] = 20000
$jqArr[20,2] = 20000
$jqArr[20,3] = 20000
$jqArr[20,4] = '$/caja 15 kilos granel'
$jqArr[20,5] = 'Perú'
$jqArr[20,6] = 1333
$jqArr[20,7] = 15
$jqArr[21,0] = 40
$jqArr[21,1] = 18000
$jqArr[21,2] = 18000
$jqArr[21,3] = 18000
$jqArr[21,4] = '$/malla 20 kilos'
$jqArr[21,5] = 'Perú'
$jqArr[21,6] = 900
$jqArr[21,7] = 20
$jqArr[22,0] = 20
$jqArr[22,1] = 20000
$jqArr[22,2] = 20000
$jqArr[22,3] = 20000
$jqArr[22,4] = '$/malla 20 kilos'
$jqArr[22,5] = 'Región de Arica y Parinacota'
$jqArr[22,6] = 1000
$jqArr[22,7] = 20
$jqArr[23,0] = 10
$jqArr[23,1] = 20000
$jqArr[23,2] = 20000
$jqArr[23,3] = 20000
$jqArr[23,4] = '$/caja 15 kilos granel'
$jqArr[23,5] = 'Región de Arica y Parinacota'
$jqArr[23,6] = 1333
$jqArr[23,7] = 15
$jqArr[24,0] = 20
$jqArr[24,1] = 20000
$jqArr[24,2] = 20000
$jqArr[24,3] = 20000
$jqArr[24,4] = '$/caja 15 kilos granel'
$jqArr[24,5] = 'Región de Arica y Parinacota'
$jqArr[24,6] = 1333
$jqArr[24,7] = 15
$jqArr[25,0] = 40
$jqArr[25,1] = 17000
$jqArr[25,2] = 18000
$jqArr[25,3] = 17500
$jqArr[25,4] = '$/malla 20 kilos'
$jqArr[25,5] = 'Perú'
$jqArr[25,6] = 875
$jqArr[25,7] = 20
$jqArr[26,0] = 45
$jqArr[26,1] = 20000
$jqArr[26,2] = 20000
$jqArr[26,3] = 20000
$jqArr[26,4] = '$/caja 15 kilos granel'
$jqArr[26,5] = 'Región de Arica y Parinacota'
$jqArr[26,6] = 1333
$jqArr[26,7] = 15
$jqArr[27,0] = 15
$jqArr[27,1] = 25000
$jqArr[27,2] = 25000
$jqArr[27,3] = 25000
$jqArr[27,4] = '$/caja 15 kilos granel'
$jqArr[27,5] = 'Perú'
$jqArr[27,6] = 1667
$jqArr[27,7] = 15
$jqArr[28,0] = 18
$jqArr[28,1] = 20000
$jqArr[28,2] = 20000
$jqArr[28,3] = 20000
$jqArr[28,4] = '$/malla 20 kilos'
$jqArr[28,5] = 'Región de Arica y Parinacota'
$jqArr[28,6] = 1000
$jqArr[28,7] = 20
$jqArr[29,0] = 15
$jqArr[29,1] = 20000
$jqArr[29,2] = 20000
$jqArr[29,3] = 20000
$jqArr[29,4] = '$/malla 20 kilos'
$jqArr[29,5] = 'Perú'
$jqArr[29,6] = 1000
$jqArr[29,7] = 20
$jqArr[30,0] = 20
$jqArr[30,1] = 20000
$jqArr[30,2] = 20000
$jqArr[30,3] = 20000
$jqArr[30,4] = '$/caja 15 kilos granel'
$jqArr[30,5] = 'Región de Arica y Parinacota'
$jqArr[30,6] = 1333
$jqArr[30,7] = 15
$jqArr[31,0] = 20
$jqArr[31,1] = 20000
$jqArr[31,2] = 20000
$jqArr[31,3] = 20000
$jqArr[31,4] = '$/caja 15 kilos granel'
$jqArr[31,5] = 'Región de Arica y Parinacota'
$jqArr[31,6] = 1333
$jqArr[31,7] = 15
$jqArr[32,0] = 40
$jqArr[32,1] = 18000
$jqArr[32,2] = 18000
$jqArr[32,3] = 18000
$jqArr[32,4] = '$/malla 20 kilos'
$jqArr[32,5] = 'Perú'
$jqArr[32,6] = 900
$jqArr[32,7] = 20
$jqArr[33,0] = 40
$jqArr[33,1] = 20000
$jqArr[33,2] = 20000
$jqArr[33,3] = 20000
$jqArr[33,4] = '$/caja 15 kilos granel'
$jqArr[33,5] = 'Perú'
$jqArr[33,6] = 1333
$jqArr[33,7] = 15
$jqArr[34,0] = 50
$jqArr[34,1] = 20000
$jqArr[34,2] = 20000
$jqArr[34,3] = 20000
$jqArr[34,4] = '$/malla 20 kilos'
$jqArr[34,5] = 'Perú'
$jqArr[34,6] = 1000
$jqArr[34,7] = 20
$jqArr[35,0] = 50
$jqArr[35,1] = 20000
$jqArr[35,2] = 20000
$jqArr[35,3] = 20000
$jqArr[35,4] = '$/malla 20 kilos'
$jqArr[35,5] = 'Perú'
$jqArr[35,6] = 1000
$jqArr[35,7] = 20
$ws.Range("J2:Q37").Value = $jqArr
